$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 41
$ws.Range("H41").Value = 1351.3334
$ws.Range("J41").Value = 1270.125
$ws.Range("L41").Value = 1270.125
$ws.Range("N41").Value = -2150.125

# ALC!row 62
$ws.Range("H62").Value = 6967.7144
$ws.Range("I62").Value = 6795.8335
$ws.Range("K62").Value = 6795.8335
$ws.Range("M62").Value = -6171.8335

# ALC!row 65
$ws.Range("H65").Value = 6967.7144
$ws.Range("I65").Value = 6795.8335
$ws.Range("K65").Value = 33979.1675
$ws.Range("M65").Value = -30859.1675

# ALC!row 111
$ws.Range("H111").Value = 11115751
$ws.Range("I111").Value = 15877145
$ws.Range("K111").Value = 47631435
$ws.Range("M111").Value = -47628368

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 32
$ws.Range("H32").Value = 8540.65
$ws.Range("I32").Value = 4534.1333
$ws.Range("J32").Value = 20560.2
$ws.Range("K32").Value = 4534.1333
$ws.Range("L32").Value = 20560.2
$ws.Range("M32").Value = -4247.1333
$ws.Range("N32").Value = -21134.2

# ARM!row 45
$ws.Range("H45").Value = 6157866
$ws.Range("I45").Value = 8549290
$ws.Range("K45").Value = 8549290
$ws.Range("M45").Value = -8548913

# ARM!row 61
$ws.Range("H61").Value = 4209.9546
$ws.Range("I61").Value = 4136.2
$ws.Range("K61").Value = 4136.2
$ws.Range("M61").Value = -3924.2

# ARM!row 74
$ws.Range("H74").Value = 23777.977
$ws.Range("I74").Value = 1245.2069
$ws.Range("K74").Value = 1245.2069
$ws.Range("M74").Value = -371.2068999999999

# ARM!row 77
$ws.Range("H77").Value = 23777.977
$ws.Range("I77").Value = 1245.2069
$ws.Range("K77").Value = 6226.0345
$ws.Range("M77").Value = -1858.0345

# ARM!row 97
$ws.Range("H97").Value = 1348673.2
$ws.Range("I97").Value = 2157345.5
$ws.Range("J97").Value = 886.2222
$ws.Range("K97").Value = 2157345.5
$ws.Range("L97").Value = 886.2222
$ws.Range("M97").Value = -2156849.5
$ws.Range("N97").Value = -1878.2222

# ARM!row 122
$ws.Range("H122").Value = 538135.75
$ws.Range("I122").Value = 2774.348
$ws.Range("J122").Value = 1307717.8
$ws.Range("K122").Value = 8323.044
$ws.Range("L122").Value = 3923153.4
$ws.Range("M122").Value = -5873.044
$ws.Range("N122").Value = -3928053.4

# ARM!row 136
$ws.Range("H136").Value = 4209.9546
$ws.Range("I136").Value = 4136.2
$ws.Range("K136").Value = 12408.6
$ws.Range("M136").Value = -9858.599999999999

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 28
$ws.Range("H28").Value = 29166.334
$ws.Range("J28").Value = 29166.334
$ws.Range("L28").Value = 29166.334
$ws.Range("N28").Value = -29656.334

# CRP!row 31
$ws.Range("H31").Value = 24174.617
$ws.Range("I31").Value = 3943.4
$ws.Range("K31").Value = 3943.4
$ws.Range("M31").Value = -3648.4

# CRP!row 34
$ws.Range("H34").Value = 24174.617
$ws.Range("I34").Value = 3943.4
$ws.Range("K34").Value = 3943.4
$ws.Range("M34").Value = -3741.4

# CRP!row 35
$ws.Range("H35").Value = 6372.222
$ws.Range("I35").Value = 2475.1667
$ws.Range("J35").Value = 14166.333
$ws.Range("K35").Value = 2475.1667
$ws.Range("L35").Value = 14166.333
$ws.Range("M35").Value = -2181.1667
$ws.Range("N35").Value = -14754.333

# CRP!row 58
$ws.Range("H58").Value = 4769.2
$ws.Range("I58").Value = 5414.067
$ws.Range("K58").Value = 5414.067
$ws.Range("M58").Value = -5211.067

# CRP!row 136
$ws.Range("H136").Value = 4769.2
$ws.Range("I136").Value = 5414.067
$ws.Range("K136").Value = 16242.201
$ws.Range("M136").Value = -13692.201

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 2
$ws.Range("H2").Value = 205.7255
$ws.Range("I2").Value = 99.878784
$ws.Range("J2").Value = 399.77777
$ws.Range("K2").Value = 599.272704
$ws.Range("L2").Value = 2398.66662
$ws.Range("M2").Value = -486.272704
$ws.Range("N2").Value = -2624.66662

# CUL!row 34
$ws.Range("H34").Value = 94.57143000000001
$ws.Range("I34").Value = 99.083336
$ws.Range("J34").Value = 67.5
$ws.Range("K34").Value = 297.250008
$ws.Range("L34").Value = 202.5
$ws.Range("M34").Value = -213.250008
$ws.Range("N34").Value = -370.5

# CUL!row 88
$ws.Range("H88").Value = 13857.143
$ws.Range("J88").Value = 13857.143
$ws.Range("L88").Value = 41571.429
$ws.Range("N88").Value = -42427.429

# CUL!row 91
$ws.Range("H91").Value = 13857.143
$ws.Range("J91").Value = 13857.143
$ws.Range("L91").Value = 41571.429
$ws.Range("N91").Value = -44535.429

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 24
$ws.Range("H24").Value = 15148.167
$ws.Range("J24").Value = 16377.8
$ws.Range("L24").Value = 16377.8
$ws.Range("N24").Value = -16723.8

# GSM!row 126
$ws.Range("H126").Value = 3906996
$ws.Range("I126").Value = 7578582.5
$ws.Range("J126").Value = 3091087.8
$ws.Range("K126").Value = 22735747.5
$ws.Range("L126").Value = 9273263.399999999
$ws.Range("M126").Value = -22733277.5
$ws.Range("N126").Value = -9278203.399999999

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 40
$ws.Range("H40").Value = 14999
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 14999
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 14999
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15271

# LTW!row 100
$ws.Range("H100").Value = 5366.593
$ws.Range("I100").Value = 5769.9
$ws.Range("J100").Value = 4214.2856
$ws.Range("K100").Value = 5769.9
$ws.Range("L100").Value = 4214.2856
$ws.Range("M100").Value = -5228.9
$ws.Range("N100").Value = -5296.2856

# LTW!row 122
$ws.Range("H122").Value = 5953.1304
$ws.Range("I122").Value = 3384.75
$ws.Range("K122").Value = 10154.25
$ws.Range("M122").Value = -7704.25

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 31
$ws.Range("H31").Value = 20003
$ws.Range("I31").Value = 18666.334
$ws.Range("K31").Value = 18666.334
$ws.Range("M31").Value = -18318.334

# WVR!row 100
$ws.Range("H100").Value = 1953.8182
$ws.Range("I100").Value = 2817
$ws.Range("J100").Value = 443.25
$ws.Range("K100").Value = 5634
$ws.Range("L100").Value = 886.5
$ws.Range("M100").Value = -5093
$ws.Range("N100").Value = -1968.5

# WVR!row 107
$ws.Range("H107").Value = 90909720
$ws.Range("I107").Value = 100000530
$ws.Range("K107").Value = 300001590
$ws.Range("M107").Value = -299999670

# WVR!row 122
$ws.Range("H122").Value = 3477.4546
$ws.Range("I122").Value = 1669.625
$ws.Range("J122").Value = 8298.333000000001
$ws.Range("K122").Value = 5008.875
$ws.Range("L122").Value = 24894.999
$ws.Range("M122").Value = -2558.875
$ws.Range("N122").Value = -29794.999

# WVR!row 126
$ws.Range("H126").Value = 1495.1613
$ws.Range("I126").Value = 1276.2106
$ws.Range("J126").Value = 1841.8334
$ws.Range("K126").Value = 3828.6318
$ws.Range("L126").Value = 5525.5002
$ws.Range("M126").Value = -1358.6318
$ws.Range("N126").Value = -10465.5002

# WVR!row 132
$ws.Range("H132").Value = 12331053
$ws.Range("I132").Value = 14928875
$ws.Range("J132").Value = 727450.25
$ws.Range("K132").Value = 44786625
$ws.Range("L132").Value = 2182350.75
$ws.Range("M132").Value = -44784095
$ws.Range("N132").Value = -2187410.75
